# Update to get data from coingecko if yahoo not working:
# refresh purchase-price/amount figures and drop the BANANA row (no longer tracked).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SOL purchase price refreshed
$ws.Range("B2").Value = 72.2

# RAY amount/price refreshed
$ws.Range("B8").Value = 4075
$ws.Range("C8").Value = 1.4

# Remove the BANANA row entirely (row 11); remaining rows shift up.
$ws.Rows.Item(11).Delete()

# After the shift, APT now sits on row 22 - refresh its amount
$ws.Range("B22").Value = 300

# After the shift, ONDO now sits on row 26 - refresh its amount/price
$ws.Range("B26").Value = 1852
$ws.Range("C26").Value = 1.134

# Restore a sensible view: scrolled near the top, B2 selected.
$ws.Range("B2").Select()
